# Updates currentAveragePrice / LevePrice / LeveProfit figures for each
# Disciple-of-the-Hand leve sheet, per the latest market-board pull from
# the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 3472
$ws.Range("I12").Value = 709.8
$ws.Range("J12").Value = 6924.75
$ws.Range("K12").Value = 709.8
$ws.Range("L12").Value = 6924.75
$ws.Range("M12").Value = -539.8
$ws.Range("N12").Value = -7264.75
# Row 40
$ws.Range("H40").Value = 1161.5
$ws.Range("I40").Value = 1096.84
$ws.Range("J40").Value = 1392.4286
$ws.Range("K40").Value = 1096.84
$ws.Range("L40").Value = 1392.4286
$ws.Range("M40").Value = -921.8399999999999
$ws.Range("N40").Value = -1742.4286
# Row 41
$ws.Range("H41").Value = 1573.4
$ws.Range("I41").Value = 2173.1428
$ws.Range("K41").Value = 2173.1428
$ws.Range("M41").Value = -1733.1428
# Row 132
$ws.Range("H132").Value = 15199002
$ws.Range("I132").Value = 15922716
$ws.Range("K132").Value = 47768148
$ws.Range("M132").Value = -47765618
# Row 137
$ws.Range("H137").Value = 1990.3334
$ws.Range("I137").Value = 1953.9286
$ws.Range("K137").Value = 5861.7858
$ws.Range("M137").Value = -3311.7858
# Row 141
$ws.Range("H141").Value = 3994.6667
$ws.Range("I141").Value = 992
$ws.Range("K141").Value = 2976
$ws.Range("M141").Value = 2204

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4314.1665
$ws.Range("I2").Value = 4899
$ws.Range("J2").Value = 1390
$ws.Range("K2").Value = 4899
$ws.Range("L2").Value = 1390
$ws.Range("M2").Value = -4786
$ws.Range("N2").Value = -1616
# Row 32
$ws.Range("H32").Value = 2065.543
$ws.Range("I32").Value = 1835.1177
$ws.Range("K32").Value = 1835.1177
$ws.Range("M32").Value = -1548.1177
# Row 74
$ws.Range("H74").Value = 12912.111
$ws.Range("I74").Value = 2035.3334
$ws.Range("K74").Value = 2035.3334
$ws.Range("M74").Value = -1161.3334
# Row 77
$ws.Range("H77").Value = 12912.111
$ws.Range("I77").Value = 2035.3334
$ws.Range("K77").Value = 10176.667
$ws.Range("M77").Value = -5808.666999999999
# Row 116
$ws.Range("H116").Value = 4314.1665
$ws.Range("I116").Value = 4899
$ws.Range("J116").Value = 1390
$ws.Range("K116").Value = 4899
$ws.Range("L116").Value = 1390
$ws.Range("M116").Value = -2605
$ws.Range("N116").Value = -5978
# Row 122
$ws.Range("H122").Value = 2043.12
$ws.Range("I122").Value = 1890.7
$ws.Range("J122").Value = 2652.8
$ws.Range("K122").Value = 5672.1
$ws.Range("L122").Value = 7958.400000000001
$ws.Range("M122").Value = -3222.1
$ws.Range("N122").Value = -12858.4
# Row 132
$ws.Range("H132").Value = 3343.8823
$ws.Range("I132").Value = 3123.0667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 9369.2001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6839.2001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4314.1665
$ws.Range("I3").Value = 4899
$ws.Range("J3").Value = 1390
$ws.Range("K3").Value = 4899
$ws.Range("L3").Value = 1390
$ws.Range("M3").Value = -4785
$ws.Range("N3").Value = -1618
# Row 86
$ws.Range("H86").Value = 1749.5
$ws.Range("J86").Value = 1666
$ws.Range("L86").Value = 1666
$ws.Range("N86").Value = -3912
# Row 89
$ws.Range("H89").Value = 1749.5
$ws.Range("J89").Value = 1666
$ws.Range("L89").Value = 8330
$ws.Range("N89").Value = -19562
# Row 105
$ws.Range("H105").Value = 3671.1428
$ws.Range("I105").Value = 3476.1765
$ws.Range("J105").Value = 4499.75
$ws.Range("K105").Value = 3476.1765
$ws.Range("L105").Value = 4499.75
$ws.Range("M105").Value = -1729.1765
$ws.Range("N105").Value = -7993.75
# Row 134
$ws.Range("H134").Value = 2756.5264
$ws.Range("I134").Value = 2471.6
$ws.Range("J134").Value = 3825
$ws.Range("K134").Value = 7414.799999999999
$ws.Range("L134").Value = 11475
$ws.Range("M134").Value = -4879.799999999999
$ws.Range("N134").Value = -16545

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3446
$ws.Range("I58").Value = 3446
$ws.Range("K58").Value = 3446
$ws.Range("M58").Value = -3243
# Row 94
$ws.Range("H94").Value = 931
$ws.Range("I94").Value = 982.4
$ws.Range("J94").Value = 802.5
$ws.Range("K94").Value = 982.4
$ws.Range("L94").Value = 802.5
$ws.Range("M94").Value = -531.4
$ws.Range("N94").Value = -1704.5
# Row 97
$ws.Range("H97").Value = 26210.25
$ws.Range("J97").Value = 26210.25
$ws.Range("L97").Value = 26210.25
$ws.Range("N97").Value = -28192.25
# Row 99
$ws.Range("H99").Value = 3911
$ws.Range("I99").Value = 3763.7693
$ws.Range("K99").Value = 3763.7693
$ws.Range("M99").Value = -2265.7693
# Row 107
$ws.Range("H107").Value = 591.17645
$ws.Range("I107").Value = 592.7143
$ws.Range("K107").Value = 592.7143
$ws.Range("M107").Value = 1327.2857
# Row 122
$ws.Range("H122").Value = 1195.5652
$ws.Range("I122").Value = 1192.579
$ws.Range("K122").Value = 3577.737
$ws.Range("M122").Value = -1127.737
# Row 126
$ws.Range("H126").Value = 3911
$ws.Range("I126").Value = 3763.7693
$ws.Range("K126").Value = 11291.3079
$ws.Range("M126").Value = -8821.3079
# Row 136
$ws.Range("H136").Value = 3446
$ws.Range("I136").Value = 3446
$ws.Range("K136").Value = 10338
$ws.Range("M136").Value = -7788

$ws = $wb.Worksheets.Item("CUL")
# Row 124
$ws.Range("H124").Value = 24785.416
$ws.Range("I124").Value = 8283.333000000001
$ws.Range("J124").Value = 27142.857
$ws.Range("K124").Value = 24849.999
$ws.Range("L124").Value = 81428.571
$ws.Range("M124").Value = -19939.999
$ws.Range("N124").Value = -91248.571
# Row 129
$ws.Range("H129").Value = 1653.6666
$ws.Range("J129").Value = 1998.25
$ws.Range("L129").Value = 5994.75
$ws.Range("N129").Value = -15994.75
# Row 131
$ws.Range("H131").Value = 52748.434
$ws.Range("J131").Value = 8127.4614
$ws.Range("L131").Value = 24382.3842
$ws.Range("N131").Value = -34462.3842

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 12653.8
$ws.Range("I70").Value = 12476.77
$ws.Range("J70").Value = 13804.5
$ws.Range("K70").Value = 12476.77
$ws.Range("L70").Value = 13804.5
$ws.Range("M70").Value = -12206.77
$ws.Range("N70").Value = -14344.5
# Row 73
$ws.Range("H73").Value = 12653.8
$ws.Range("I73").Value = 12476.77
$ws.Range("J73").Value = 13804.5
$ws.Range("K73").Value = 12476.77
$ws.Range("L73").Value = 13804.5
$ws.Range("M73").Value = -11540.77
$ws.Range("N73").Value = -15676.5
# Row 97
$ws.Range("H97").Value = 740.4375
$ws.Range("I97").Value = 648.5454999999999
$ws.Range("J97").Value = 942.6
$ws.Range("K97").Value = 648.5454999999999
$ws.Range("L97").Value = 942.6
$ws.Range("M97").Value = -152.5454999999999
$ws.Range("N97").Value = -1934.6
# Row 102
$ws.Range("H102").Value = 2384.5186
$ws.Range("I102").Value = 2473.2173
$ws.Range("J102").Value = 1874.5
$ws.Range("K102").Value = 2473.2173
$ws.Range("L102").Value = 1874.5
$ws.Range("M102").Value = -851.2172999999998
$ws.Range("N102").Value = -5118.5
# Row 107
$ws.Range("H107").Value = 1332.4706
$ws.Range("I107").Value = 1047.2
$ws.Range("K107").Value = 1047.2
$ws.Range("M107").Value = 872.8
# Row 122
$ws.Range("H122").Value = 2209.6667
$ws.Range("I122").Value = 2083.3333
$ws.Range("K122").Value = 6249.999899999999
$ws.Range("M122").Value = -3799.999899999999
# Row 126
$ws.Range("H126").Value = 30132.584
$ws.Range("I126").Value = 42949
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 128847
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -126377
$ws.Range("N126").Value = -18439.25
# Row 132
$ws.Range("H132").Value = 253040.12
$ws.Range("I132").Value = 259500.12
$ws.Range("K132").Value = 778500.36
$ws.Range("M132").Value = -775970.36

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3900.2703
$ws.Range("I16").Value = 1422.75
$ws.Range("J16").Value = 6815
$ws.Range("K16").Value = 1422.75
$ws.Range("L16").Value = 6815
$ws.Range("M16").Value = -1252.75
$ws.Range("N16").Value = -7155
# Row 46
$ws.Range("H46").Value = 1662.2727
$ws.Range("I46").Value = 1855
$ws.Range("J46").Value = 795
$ws.Range("K46").Value = 1855
$ws.Range("L46").Value = 795
$ws.Range("M46").Value = -1667
$ws.Range("N46").Value = -1171
# Row 55
$ws.Range("H55").Value = 129.17857
$ws.Range("J55").Value = 154.81818
$ws.Range("L55").Value = 154.81818
$ws.Range("N55").Value = -500.81818
# Row 122
$ws.Range("H122").Value = 7033.7896
$ws.Range("I122").Value = 6710.923
$ws.Range("K122").Value = 20132.769
$ws.Range("M122").Value = -17682.769

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1989.7693
$ws.Range("J107").Value = 3374.1667
$ws.Range("L107").Value = 10122.5001
$ws.Range("N107").Value = -13962.5001
# Row 122
$ws.Range("H122").Value = 3932.4
$ws.Range("I122").Value = 3888.7
$ws.Range("J122").Value = 4019.8
$ws.Range("K122").Value = 11666.1
$ws.Range("L122").Value = 12059.4
$ws.Range("M122").Value = -9216.099999999999
$ws.Range("N122").Value = -16959.4
# Row 126
$ws.Range("H126").Value = 4668.1055
$ws.Range("I126").Value = 4538.375
$ws.Range("K126").Value = 13615.125
$ws.Range("M126").Value = -11145.125
# Row 132
$ws.Range("H132").Value = 2226.6086
$ws.Range("I132").Value = 2282.3635
$ws.Range("K132").Value = 6847.0905
$ws.Range("M132").Value = -4317.0905
# Row 136
$ws.Range("H136").Value = 3287.5
$ws.Range("I136").Value = 3383.4443
$ws.Range("K136").Value = 10150.3329
$ws.Range("M136").Value = -7600.332900000001
# Row 139
$ws.Range("H139").Value = 52947.145
$ws.Range("J139").Value = 52947.145
$ws.Range("L139").Value = 52947.145
$ws.Range("N139").Value = -63227.145
